$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UC unit boilers")

# --- Row 20: blank spacer row styled like the other "closing" blank rows (format only, from G12 which uses style 4) ---
$ws.Range("G12").Copy()
$ws.Range("A20:J20").PasteSpecial(-4122)
$ws.Range("L20").PasteSpecial(-4122)

# --- Row 22: section header "~UC_T" (same value+format as G6) ---
$ws.Range("G6").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("G22").Value2 = $ws.Range("G6").Value2

# --- Row 23: column headers (same values+format as row 7, columns A:J) ---
$ws.Range("A7:J7").Copy()
$ws.Range("A23:J23").PasteSpecial(-4122)
$ws.Range("A23:J23").Value2 = $ws.Range("A7:J7").Value2

# --- Row 24: new UC data row ---
# F:J match row 8's F:J exactly (2019, FX, 1, 0, 5) -> copy value+format
$ws.Range("F8:J8").Copy()
$ws.Range("F24:J24").PasteSpecial(-4122)
$ws.Range("F24:J24").Value2 = $ws.Range("F8:J8").Value2
# A24, B24, D24 need the same format as row 8 (style 2) but new text
$ws.Range("A8").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("A24").Value = "UC_RSD_Cookers"
$ws.Range("B24").Value = "Number of Cookers to equal number of Houses"
$ws.Range("D24").Value = "RSDCK"

# --- Row 25: second line of the new UC ---
# H25 matches row 9's H9 exactly (-1) -> copy value+format
$ws.Range("H9").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("H25").Value2 = $ws.Range("H9").Value2
# C25 needs the same format as row 9's C9 (style 2) but a literal new value
$ws.Range("C9").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = "R-RSDCK*"

$excel.CutCopyMode = 0

# --- Update selection to match authored state ---
[void]$ws.Range("C27").Select()
